# Workbook already open
$wb = $excel.ActiveWorkbook

# --- 1. Rename the 4 existing sheets ---------------------------------------
$wb.Worksheets.Item(1).Name = "PocHistone RLFP 004"
$wb.Worksheets.Item(2).Name = "PocHistone RLFP 005"
$wb.Worksheets.Item(3).Name = "PocHistone RLFP 006"
$wb.Worksheets.Item(4).Name = "PocHistone RLFP 007"

# --- 2. Update the selection on "PocHistone RLFP 007" (was tab-selected) ---
$ws007 = $wb.Worksheets.Item(4)
$ws007.Range("B1:M1").Select()

# --- 3. Add 16 new plate sheets, in order, named 008..023 -------------------
$newNames = @(
    "PocHistone RLFP 008",
    "PocHistone RLFP 009",
    "PocHistone RLFP 010",
    "PocHistone RLFP 011",
    "PocHistone RLFP 012",
    "PocHistone RLFP 013",
    "PocHistone RLFP 014",
    "PocHistone RLFP 015",
    "PocHistone RLFP 016",
    "PocHistone RLFP 017",
    "PocHistone RLFP 018",
    "PocHistone RLFP 019",
    "PocHistone RLFP 020",
    "PocHistone RLFP 021",
    "PocHistone RLFP 022",
    "PocHistone RLFP 023"
)

# Sheets "017" and "018" (indices 9 and 10 in $newNames) stay completely blank
$blankIndexes = @(9, 10)

$rowLabels = @("A","B","C","D","E","F","G","H")

$newSheets = @{}

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $last)
    $ws.Name = $newNames[$i]
    $newSheets[$i] = $ws

    if ($blankIndexes -notcontains $i) {
        # Header row: 1..12 across B1:M1
        for ($c = 2; $c -le 13; $c++) {
            $ws.Cells.Item(1, $c).Value = $c - 1
        }
        # Row labels A2:A9
        for ($r = 2; $r -le 9; $r++) {
            $ws.Cells.Item($r, 1).Value = $rowLabels[$r - 2]
        }
        # Selection covering the whole plate grid
        $ws.Range("A1:M10").Select()
    }
}

# --- 4. Final active sheet/selection: "PocHistone RLFP 008" cell B2 --------
$ws008 = $newSheets[0]
$ws008.Activate()
$ws008.Range("B2").Select()
